# "Generate Report for Handoff"
# A new handoff run produced a fresh id (898bd444-...-> 03ba9452-...) and a
# fresh target-file hash (f179d386...-> a3707126...), so every cell / hyperlink
# that echoes the old id or old handoff timestamps needs to move to the new
# ones. Only the *display* text + cell values change - the underlying
# hyperlink target addresses are left exactly as they were.

$wb = $excel.ActiveWorkbook

$oldId = "898bd444-37f3-413a-a8fe-2774988fc6b7"
$newId = "03ba9452-c133-46a9-a560-bb223b0e6fad"

$oldHash = "f179d3862d59529bdd7fec6240f701088dd1499d"
$newHash = "a3707126b6a17ea28492e81614a809732deaf9b4"

# ---------------------------------------------------------------------------
# Sheet "Overview": only one hyperlink (A2 -> the .md handoff file)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$newMdName = "$newId.md"
$aAddr = $ws.Hyperlinks.Item(1).Address

$ws.Range("A2").Value = $newMdName
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $aAddr, "", "", $newMdName)
$ws.Range("A2").Style = "HyperLink"

$ws.Range("D2").Value = "2016-48-17 18:48:53"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": hyperlinks on A2 (.md), B2 (.md extension, unchanged) and
# D2 (the zh-cn .xlf handoff file)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$aAddr = $ws.Hyperlinks.Item(1).Address
$bAddr = $ws.Hyperlinks.Item(2).Address
$bDisp = $ws.Hyperlinks.Item(2).TextToDisplay
$dAddr = $ws.Hyperlinks.Item(3).Address

$newXlfNameZh = "$newId.$newHash.zh-cn.xlf"

$ws.Range("A2").Value = $newMdName
$ws.Range("D2").Value = $newXlfNameZh

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $aAddr, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("B2"), $bAddr, "", "", $bDisp)
$ws.Hyperlinks.Add($ws.Range("D2"), $dAddr, "", "", $newXlfNameZh)
$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Style = "HyperLink"
$ws.Range("D2").Style = "HyperLink"

$ws.Range("E2").Value = "2016-03-17 18:48:49"

# ---------------------------------------------------------------------------
# Sheet "de-de": hyperlinks on A2 (.md), B2 (.md extension, unchanged) and
# D2 (the de-de .xlf handoff file)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$aAddr = $ws.Hyperlinks.Item(1).Address
$bAddr = $ws.Hyperlinks.Item(2).Address
$bDisp = $ws.Hyperlinks.Item(2).TextToDisplay
$dAddr = $ws.Hyperlinks.Item(3).Address

$newXlfNameDe = "$newId.$newHash.de-de.xlf"

$ws.Range("A2").Value = $newMdName
$ws.Range("D2").Value = $newXlfNameDe

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $aAddr, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("B2"), $bAddr, "", "", $bDisp)
$ws.Hyperlinks.Add($ws.Range("D2"), $dAddr, "", "", $newXlfNameDe)
$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Style = "HyperLink"
$ws.Range("D2").Style = "HyperLink"

$ws.Range("E2").Value = "2016-03-17 18:48:53"
